# "Generate Report for Handoff"
#
# The localization pipeline re-ran and produced a new handoff package for
# e2e\b.md (b.63290e5768f688058c7b37413b0a5c26c308f864.*.xlf). Reflect that
# on all three sheets:
#   - Overview: row for b.md moves from "Handed back: in sync with en-US"
#     to "Ready for handoff", with the new generation timestamp.
#   - zh-cn / de-de detail sheets: b.md's Status, Content Duplicate flag,
#     Latest Handoff File, Latest Handoff Datetime and Error Detail are
#     updated to match; the Error Detail column is widened to fit the new,
#     much longer message.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4993103581bd6e6f1c51e95a099fb29337e72487/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d1ad0ed080983e4062a36685b9847a69d738011c/e2e/b.md."

# ---------------------------------------------------------------------
# Overview sheet: row 3 is e2e\b.md
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-11-15 16:46:09"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is b.md
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-11-15 16:45:53"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Range("P1:P3").ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet: row 3 is b.md
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-11-15 16:46:09"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Range("P1:P3").ColumnWidth = 39.17
